# Auto-generated edit script: updates crypto price/volume table
# to reflect refreshed data from the GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.193.94"
$ws.Range("E2").Value = "'  -4.62%  "
$ws.Range("D3").Value = "'2.481.74"
$ws.Range("E3").Value = "'  -4.23%  "
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("D5").Value = "'538.14"
$ws.Range("E5").Value = "'  -2.85%  "
$ws.Range("D6").Value = "'145.70"
$ws.Range("E6").Value = "'  -6.21%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "'  -0.18%  "
$ws.Range("E8").Value = "'  -3.53%  "
$ws.Range("D9").Value = "'2.511.27"
$ws.Range("E9").Value = "'  -3.41%  "
$ws.Range("E10").Value = "'  -3.84%  "
$ws.Range("E11").Value = "'  -1.98%  "
$ws.Range("D12").Value = "'5.61"
$ws.Range("E12").Value = "'  +2.11%  "
$ws.Range("D13").Value = "'0.357"
$ws.Range("E13").Value = "'  -2.34%  "
$ws.Range("D14").Value = "'2.927.59"
$ws.Range("E14").Value = "'  -3.95%  "
$ws.Range("D15").Value = "'24.14"
$ws.Range("E15").Value = "'  -5.83%  "
$ws.Range("D16").Value = "'59.108.32"
$ws.Range("E16").Value = "'  -4.57%  "
$ws.Range("D17").Value = "'0.0000138"
$ws.Range("E17").Value = "'  -3.77%  "
$ws.Range("D18").Value = "'2.511.47"
$ws.Range("E18").Value = "'  -3.13%  "
$ws.Range("D19").Value = "'11.48"
$ws.Range("E19").Value = "'  -1.65%  "
$ws.Range("E20").Value = "'  -4.61%  "
$ws.Range("D21").Value = "'326.07"
$ws.Range("E21").Value = "'  -3.99%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "'  -0.02%  "
$ws.Range("E23").Value = "'  -4.31%  "
$ws.Range("D24").Value = "'61.13"
$ws.Range("E24").Value = "'  -2.32%  "
$ws.Range("D25").Value = "'0.444"
$ws.Range("E25").Value = "'  -11.10%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.162"
$ws.Range("E26").Value = "'  -3.76%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "'2.617.31"
$ws.Range("E27").Value = "'  -3.39%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'0.995"
$ws.Range("E28").Value = "'  -0.49%  "
$ws.Range("E29").Value = "'  -2.76%  "
$ws.Range("D30").Value = "'7.09"
$ws.Range("E30").Value = "'  -0.35%  "
$ws.Range("D31").Value = "'1.29"
$ws.Range("E31").Value = "'  -2.02%  "
$ws.Range("D32").Value = "'0.0₃0782"
$ws.Range("E32").Value = "'  -6.88%  "
$ws.Range("E33").Value = "'  -5.77%  "
$ws.Range("D34").Value = "'0.996"
$ws.Range("E34").Value = "'  -0.21%  "
$ws.Range("D35").Value = "'158.30"
$ws.Range("E35").Value = "'  -1.12%  "
$ws.Range("E36").Value = "'  -0.59%  "
$ws.Range("D37").Value = "'18.61"
$ws.Range("E37").Value = "'  -3.43%  "
$ws.Range("E38").Value = "'  -5.43%  "
$ws.Range("E39").Value = "'  -7.31%  "
$ws.Range("D40").Value = "'5.88"
$ws.Range("E40").Value = "'  -2.79%  "
$ws.Range("D41").Value = "'313.75"
$ws.Range("E41").Value = "'  -8.18%  "
$ws.Range("D42").Value = "'36.79"
$ws.Range("E42").Value = "'  -2.06%  "
$ws.Range("D43").Value = "'3.74"
$ws.Range("E43").Value = "'  -4.76%  "
$ws.Range("D44").Value = "'0.830"
$ws.Range("E44").Value = "'  -7.74%  "
$ws.Range("D45").Value = "'0.996"
$ws.Range("E45").Value = "'  -0.18%  "
$ws.Range("E46").Value = "'  -1.73%  "
$ws.Range("D47").Value = "'10.75"
$ws.Range("E47").Value = "'  -1.74%  "
$ws.Range("D48").Value = "'125.20"
$ws.Range("E48").Value = "'  -0.25%  "
$ws.Range("D49").Value = "'0.0526"
$ws.Range("E49").Value = "'  -4.22%  "
$ws.Range("D50").Value = "'0.0931"
$ws.Range("E50").Value = "'  -3.64%  "
$ws.Range("D51").Value = "'0.0230"
$ws.Range("E51").Value = "'  -4.30%  "
